$wb = $excel.ActiveWorkbook

# Rename the first sheet "Dati generali" -> "DatiGenerali"
$wsDati = $wb.Worksheets.Item(1)
$wsDati.Name = "DatiGenerali"

$wsMisurazioni = $wb.Worksheets.Item(2)
$wsGrafico = $wb.Worksheets.Item(3)

# Update the selection on "Misurazioni" (no longer the active tab)
$wsMisurazioni.Range("F16").Select()

# Update the selection on "Grafico" (no longer the active tab, so it stops
# rendering/recalculating its embedded chart until the workbook is reopened
# and that sheet is revisited)
$wsGrafico.Range("O15").Select()

# Make "DatiGenerali" the active sheet/tab with a fresh selection at B1
$wsDati.Activate()
$wsDati.Range("B1").Select()
